# Regenerate orders with updated distance/sizes.
# The experiment's distance and size codes were renumbered:
#   D80 -> D86, D64 -> D69, D51 -> D55, S30 -> S31
# (S25 and S20 are left unchanged). This affects every text column that
# embeds a distance or size code: Condition (B), Filename_Left (D),
# Filename_Right (E), Distance (H) and Size (J).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowCount = $ws.UsedRange.Rows.Count
$cols = @(2, 4, 5, 8, 10)   # B, D, E, H, J

for ($r = 2; $r -le $rowCount; $r++) {
    foreach ($c in $cols) {
        $cell = $ws.Cells.Item($r, $c)
        $v = $cell.Value()
        if ($v -ne $null) {
            $nv = $v.Replace("D80", "D86").Replace("D64", "D69").Replace("D51", "D55").Replace("S30", "S31")
            if ($nv -ne $v) {
                $cell.Value = $nv
            }
        }
    }
}
